$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column Q (string: themeTrack) ---
$ws.Range("Q4").Value = "string"
$ws.Range("Q5").Value = "themeTrack"
$ws.Range("Q6").Value = "ui/assets/topdown/top-down-shooter/music/theme-1.ogg"
$ws.Range("Q7").Value = "ui/assets/topdown/top-down-shooter/music/theme-2.ogg"
$ws.Range("Q8").Value = "ui/assets/topdown/top-down-shooter/music/theme-3.ogg"

# --- Row 6 (operator 0001): localized name/role + new weapon/relic ids ---
$ws.Range("D6").Value = "萨布尔「霓虹猎手」"
$ws.Range("E6").Value = "霓虹清剿者"
$ws.Range("F6").Value = "weapon:20020001"
$ws.Range("G6").Value = "relic:30050001"

# --- Row 7 (operator 0002): localized name/role + new weapon/relic ids ---
$ws.Range("D7").Value = "伊莉丝「破晓灯塔」"
$ws.Range("E7").Value = "共鸣信标师"
$ws.Range("F7").Value = "weapon:20020002"
$ws.Range("G7").Value = "relic:30050002"

# --- Row 8 (operator 0003): localized name/role + new weapon/relic ids ---
$ws.Range("D8").Value = "马洛「潮汐制裁」"
$ws.Range("E8").Value = "深渊猎手"
$ws.Range("F8").Value = "weapon:20020003"
$ws.Range("G8").Value = "relic:30050003"

# spriteScale for row 8 changed 0.95 -> 0.88; this column is stored as
# text (numberStoredAsText), so keep it text by using Excel's leading
# apostrophe text-prefix instead of letting it coerce to a number.
$ws.Range("P8").Value = "'0.88"
